{"js": "const pairs = [\n  [\"82\u00d738=3116\", \"17\u00d749=833\"],\n  [\"76\u00d767=5092\", \"11\u00d796=1056\"],\n  [\"68\u00d739=2652\", \"62\u00d732=1984\"],\n  [\"69\u00d790=6210\", \"94\u00d741=3854\"],\n  [\"59\u00d717=1003\", \"67\u00d761=4087\"],\n  [\"61\u00d765=3965\", \"87\u00d728=2436\"],\n  [\"27\u00d789=2403\", \"83\u00d711=913\"],\n  [\"19\u00d797=1843\", \"42\u00d716=672\"],\n  [\"82\u00d751=4182\", \"67\u00d718=1206\"],\n  [\"47\u00d761=2867\", \"69\u00d724=1656\"],\n  [\"39\u00d766=2574\", \"83\u00d763=5229\"],\n  [\"43\u00d753=2279\", \"24\u00d741=984\"],\n  [\"98\u00d750=4900\", \"47\u00d793=4371\"],\n  [\"24\u00d720=480\", \"51\u00d774=3774\"],\n  [\"91\u00d760=5460\", \"28\u00d750=1400\"],\n  [\"26\u00d775=1950\", \"54\u00d753=2862\"],\n  [\"23\u00d736=828\", \"58\u00d746=2668\"],\n  [\"13\u00d775=975\", \"34\u00d714=476\"],\n  [\"69\u00d789=6141\", \"96\u00d723=2208\"],\n  [\"72\u00d764=4608\", \"96\u00d764=6144\"],\n  [\"43\u00d760=2580\", \"95\u00d760=5700\"],\n  [\"84\u00d749=4116\", \"78\u00d725=1950\"],\n  [\"82\u00d752=4264\", \"92\u00d722=2024\"],\n  [\"13\u00d733=429\", \"36\u00d790=3240\"],\n  [\"52\u00d775=3900\", \"27\u00d780=2160\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '82\u00d738=3116'\n$find.Replacement.Text = '17\u00d749=833'\n$find.Execute([ref]'82\u00d738=3116', $false, $false, $false, $false, $false, $true, 1, $false, '17\u00d749=833', 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '76\u00d767=5092'\n$find.Replacement.Text = '11\u00d796=1056'\n$find.Execute([ref]'76\u00d767=5092', $false, $false, $false, $false, $false, $true, 1, $false, '11\u00d796=1056', 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '68\u00d739=2652'\n$find.Replacement.Text = '62\u00d732=1984'\n$find.Execute([ref]'68\u00d739=2652', $false, $false, $false, $false, $false, $true, 1, $false, '62\u00d732=1984', 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '69\u00d790=6210'\n$find.Replacement.Text = '94\u00d741=3854'\n$find.Execute([ref]'69\u00d790=6210', $false, $false, $false, $false, $false, $true, 1, $false, '94\u00d741=3854', 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '59\u00d717=1003'\n$find.Replacement.Text = '67\u00d761=4087'\n$find.Execute([ref]'59\u00d717=1003', $false, $false, $false, $false, $false, $true, 1, $false, '67\u00d761=4087', 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '61\u00d765=3965'\n$find.Replacement.Text = '87\u00d728=2436'\n$find.Execute([ref]'61\u00d765=3965', $false, $false, $false, $false, $false, $true, 1, $false, '87\u00d728=2436', 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '27\u00d789=2403'\n$find.Replacement.Text = '83\u00d711=913'\n$find.Execute([ref]'27\u00d789=2403', $false, $false, $false, $false, $false, $true, 1, $false, '83\u00d711=913', 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '19\u00d797=1843'\n$find.Replacement.Text = '42\u00d716=672'\n$find.Execute([ref]'19\u00d797=1843', $false, $false, $false, $false, $false, $true, 1, $false, '42\u00d716=672', 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '82\u00d751=4182'\n$find.Replacement.Text = '67\u00d718=1206'\n$find.Execute([ref]'82\u00d751=4182', $false, $false, $false, $false, $false, $true, 1, $false, '67\u00d718=1206', 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '47\u00d761=2867'\n$find.Replacement.Text = '69\u00d724=1656'\n$find.Execute([ref]'47\u00d761=2867', $false, $false, $false, $false, $false, $true, 1, $false, '69\u00d724=1656', 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '39\u00d766=2574'\n$find.Replacement.Text = '83\u00d763=5229'\n$find.Execute([ref]'39\u00d766=2574', $false, $false, $false, $false, $false, $true, 1, $false, '83\u00d763=5229', 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '43\u00d753=2279'\n$find.Replacement.Text = '24\u00d741=984'\n$find.Execute([ref]'43\u00d753=2279', $false, $false, $false, $false, $false, $true, 1, $false, '24\u00d741=984', 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '98\u00d750=4900'\n$find.Replacement.Text = '47\u00d793=4371'\n$find.Execute([ref]'98\u00d750=4900', $false, $false, $false, $false, $false, $true, 1, $false, '47\u00d793=4371', 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '24\u00d720=480'\n$find.Replacement.Text = '51\u00d774=3774'\n$find.Execute([ref]'24\u00d720=480', $false, $false, $false, $false, $false, $true, 1, $false, '51\u00d774=3774', 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '91\u00d760=5460'\n$find.Replacement.Text = '28\u00d750=1400'\n$find.Execute([ref]'91\u00d760=5460', $false, $false, $false, $false, $false, $true, 1, $false, '28\u00d750=1400', 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '26\u00d775=1950'\n$find.Replacement.Text = '54\u00d753=2862'\n$find.Execute([ref]'26\u00d775=1950', $false, $false, $false, $false, $false, $true, 1, $false, '54\u00d753=2862', 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '23\u00d736=828'\n$find.Replacement.Text = '58\u00d746=2668'\n$find.Execute([ref]'23\u00d736=828', $false, $false, $false, $false, $false, $true, 1, $false, '58\u00d746=2668', 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '13\u00d775=975'\n$find.Replacement.Text = '34\u00d714=476'\n$find.Execute([ref]'13\u00d775=975', $false, $false, $false, $false, $false, $true, 1, $false, '34\u00d714=476', 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '69\u00d789=6141'\n$find.Replacement.Text = '96\u00d723=2208'\n$find.Execute([ref]'69\u00d789=6141', $false, $false, $false, $false, $false, $true, 1, $false, '96\u00d723=2208', 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '72\u00d764=4608'\n$find.Replacement.Text = '96\u00d764=6144'\n$find.Execute([ref]'72\u00d764=4608', $false, $false, $false, $false, $false, $true, 1, $false, '96\u00d764=6144', 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '43\u00d760=2580'\n$find.Replacement.Text = '95\u00d760=5700'\n$find.Execute([ref]'43\u00d760=2580', $false, $false, $false, $false, $false, $true, 1, $false, '95\u00d760=5700', 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '84\u00d749=4116'\n$find.Replacement.Text = '78\u00d725=1950'\n$find.Execute([ref]'84\u00d749=4116', $false, $false, $false, $false, $false, $true, 1, $false, '78\u00d725=1950', 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '82\u00d752=4264'\n$find.Replacement.Text = '92\u00d722=2024'\n$find.Execute([ref]'82\u00d752=4264', $false, $false, $false, $false, $false, $true, 1, $false, '92\u00d722=2024', 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '13\u00d733=429'\n$find.Replacement.Text = '36\u00d790=3240'\n$find.Execute([ref]'13\u00d733=429', $false, $false, $false, $false, $false, $true, 1, $false, '36\u00d790=3240', 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '52\u00d775=3900'\n$find.Replacement.Text = '27\u00d780=2160'\n$find.Execute([ref]'52\u00d775=3900', $false, $false, $false, $false, $false, $true, 1, $false, '27\u00d780=2160', 2) | Out-Null\n"}
